{"js": "// Apply the \"Added many more features\" edit to the Dim Sum Prize review.\n//\n// Each entry is an exact-text replacement. The title/heading string\n// (\"Play Dim Sum Prize for Free - Delicious Chinese-themed Slot\")\n// occurs twice in the document (the H1 heading and the bold recap line\n// near the end) so we replace every match found by search(), not just\n// the first one.\nconst replacements = [\n  [\n    \"Play Dim Sum Prize for Free - Delicious Chinese-themed Slot\",\n    \"Play Dim Sum Prize Free - Exciting Slot Game Review\",\n  ],\n  [\n    \"High RTP of 97.18%, great for experienced players\",\n    \"High RTP of 97.18%\",\n  ],\n  [\n    \"Beautiful graphics that bring the Chinese restaurant theme to life\",\n    \"Impressive graphics and design\",\n  ],\n  [\n    \"Betting range suitable for both low- and high-stakes players\",\n    \"Wide betting range for low- and high-stakes players\",\n  ],\n  [\n    \"Autoplay feature allows for uninterrupted gameplay\",\n    \"Two exciting bonus features\",\n  ],\n  [\n    \"Limited number of paylines may not appeal to some players\",\n    \"Limited number of paylines (10)\",\n  ],\n  [\n    \"Free spins can be difficult to trigger\",\n    \"Free spins feature can be hard to trigger\",\n  ],\n  [\n    \"Dim Sum Prize is a Chinese-themed slot game with 10 fixed paylines. Play now for free and enjoy two bonus features and an impressive RTP of 97.18%.\",\n    \"Play Dim Sum Prize for free and enjoy impressive graphics, exciting bonus features, and a wide betting range.\",\n  ],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the \"Added many more features\" edit to the Dim Sum Prize review.\n#\n# Each pair is an exact whole-string Find & Replace. The title/heading\n# string (\"Play Dim Sum Prize for Free - Delicious Chinese-themed Slot\")\n# occurs twice in the document (the H1 heading and the bold recap line\n# near the end); Replace:=wdReplaceAll (2) replaces every occurrence in\n# one call, so it's fine to run it the same way as the other (single-\n# occurrence) pairs below.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"Play Dim Sum Prize for Free - Delicious Chinese-themed Slot\", \"Play Dim Sum Prize Free - Exciting Slot Game Review\"),\n  @(\"High RTP of 97.18%, great for experienced players\", \"High RTP of 97.18%\"),\n  @(\"Beautiful graphics that bring the Chinese restaurant theme to life\", \"Impressive graphics and design\"),\n  @(\"Betting range suitable for both low- and high-stakes players\", \"Wide betting range for low- and high-stakes players\"),\n  @(\"Autoplay feature allows for uninterrupted gameplay\", \"Two exciting bonus features\"),\n  @(\"Limited number of paylines may not appeal to some players\", \"Limited number of paylines (10)\"),\n  @(\"Free spins can be difficult to trigger\", \"Free spins feature can be hard to trigger\"),\n  @(\"Dim Sum Prize is a Chinese-themed slot game with 10 fixed paylines. Play now for free and enjoy two bonus features and an impressive RTP of 97.18%.\", \"Play Dim Sum Prize for free and enjoy impressive graphics, exciting bonus features, and a wide betting range.\")\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Execute(\n    $oldText,\n    $true,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    $newText,\n    2\n  )\n}\n"}
